# Change IRLZ44N for 2N7002
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q1) - Mosfet: swap part number, package and LCSC code
$ws.Range("B3").Value = "2N7002"
$ws.Range("C3").Value = "SOT-23-3"
$ws.Range("E3").Value = "C8545"

# Reflect the active selection after the edit
$ws.Range("E3").Select()

$wb.Save()
